$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: rename "lastname" column to "Email" ---
$ws.Range("C1").Value = "Email"

# --- Rename testcase2 first name from "Test" to "Test Peter" ---
$ws.Range("B3").Value = "Test Peter"

# --- Replace column C values with email addresses + mailto hyperlinks ---
# (Hyperlinks.Add also sets the cell's display text / value for us)
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:ASingh@mailinator.com", [System.Type]::Missing, "mailto:ASingh@mailinator.com", "ASingh@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:testpeter@mailinator.com", [System.Type]::Missing, "mailto:testpeter@mailinator.com", "testpeter@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:parker@mailinator.com", [System.Type]::Missing, [System.Type]::Missing, "parker@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:peter@mailinator.com", [System.Type]::Missing, [System.Type]::Missing, "peter@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:testm@mailinator.com", [System.Type]::Missing, [System.Type]::Missing, "testm@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:singh@mailinator.com", [System.Type]::Missing, [System.Type]::Missing, "singh@mailinator.com")

# --- Widen column C so the email addresses are fully visible ---
$ws.Columns("C").ColumnWidth = 26

# --- Move / restore the active selection to C11 ---
$ws.Range("C11").Select() | Out-Null
